# Implemented multithreading. Result is still too low.
# - Add a new worksheet "Sheet1" right after "Asset 1"
# - Fill in PARAMS / MOCKED_RANDOM_NUMBERS rows on the new sheet
# - Update active cell selections on both sheets, new sheet becomes active tab

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert a brand-new worksheet right after "Asset 1" (becomes sheet2.xml / "Sheet1")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)

# Move selection on the original sheet (no longer the active tab afterwards)
$ws1.Range("E43").Select()

# Populate the new sheet's data
$newSheet.Range("B2").Value = "PARAMS:"
$newSheet.Range("B3").Value = "MOCKED_RANDOM_NUMBERS:"
$newSheet.Range("F3").Value = 1.1
$newSheet.Range("G3").Value = 0.9
$newSheet.Range("H3").Value = 1
$newSheet.Range("I3").Value = 0.8
$newSheet.Range("J3").Value = 1.2

# Select the cell on the new sheet, making it the active tab/sheet
$newSheet.Range("O24").Select()
